$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (border style used by existing data rows) down to the
# two new rows before writing their values, mirroring how a user would
# drag-fill / copy an existing row when adding new entries.
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G7").PasteSpecial(-4122)

$data = @(
    @("Nano",  "Prototipo",            "Titulo para mi proyecto de prubas Gerry",    "Gerry Deustúa Hernández", "BOLDBGOSDBGOSANGSIGNSOPGSDG", "BOLDBGOSDBGOSANGSIGNSOPGSDG", "Titulo para mi proyecto de prubas Gerry"),
    @("Nano",  "Concepto",             "Test89",                                      "Mikel Edel",              "rrrrrrrrrr",                   "aaaaaaa",                      "Test89"),
    @("Nano",  "Prototipo finalizado", "Proyecto de prueba",                          "Marlon Martínez",         "link",                          "link",                         "Proyecto de prueba"),
    @("Nexus", "Prototipo finalizado", "fsfa",                                        "Marlon Martínez",         "fasfsa",                        "fsa",                          "fsfa"),
    @("Nexus", "Concepto",             "Robot automata para automatizar automatas",  "Gerry Deustúa Hernández", "sdgasdgasdg",                   "sadgsadg",                     "Robot automata para automatizar automatas"),
    @("Cyber", "Concepto",             "dsfomo´ghdsrg",                              "Marlon Martínez",         "dgfgj sdlgondskgdsfg",          "df´pdsfkg´pdgksdg",            "dsfomo´ghdsrg")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
